# 14-Apr-2024: Administrator functions were implemented.
# Adds two new "administrator" rows (candidates data / test results) above
# the existing "test paper" / "marksheet" rows on Sheet1, and moves the
# selection to the new first inserted cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing rows 8-9 ("test paper"/"marksheet") down to 10-11 and
# open up two blank rows at 8-9, inheriting the surrounding formatting.
[void]$ws.Rows("8:9").Insert()

# Fill in the two new rows. Values are entered in the same order the
# original author typed them (B8, B9, A9, A8) so the shared-string table
# is built up in the same sequence.
$ws.Range("B8").Value = "candidates.xlsx"
$ws.Range("B9").Value = "scores.xlsx"
$ws.Range("A9").Value = "test results"
$ws.Range("A8").Value = "candidates data"

# Keep the new rows the same height as the rest of the table.
$ws.Rows("8:9").RowHeight = 19

# Move the active selection to the first new row, matching the saved view.
[void]$ws.Range("A8").Select()
